$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1927.3077
$ws.Range("I40").Value = 2147.4736
$ws.Range("J40").Value = 1329.7142
$ws.Range("K40").Value = 2147.4736
$ws.Range("L40").Value = 1329.7142
$ws.Range("M40").Value = -1972.4736
$ws.Range("N40").Value = -1679.7142

$ws.Range("H76").Value = 3521.2632
$ws.Range("I76").Value = 3326.6667
$ws.Range("J76").Value = 4251
$ws.Range("K76").Value = 3326.6667
$ws.Range("L76").Value = 4251
$ws.Range("M76").Value = -3011.6667
$ws.Range("N76").Value = -4881

$ws.Range("H79").Value = 3521.2632
$ws.Range("I79").Value = 3326.6667
$ws.Range("J79").Value = 4251
$ws.Range("K79").Value = 3326.6667
$ws.Range("L79").Value = 4251
$ws.Range("M79").Value = -2234.6667
$ws.Range("N79").Value = -6435

$ws.Range("H86").Value = 45487.39
$ws.Range("I86").Value = 126413.625
$ws.Range("J86").Value = 2326.7334
$ws.Range("K86").Value = 126413.625
$ws.Range("L86").Value = 2326.7334
$ws.Range("M86").Value = -125290.625
$ws.Range("N86").Value = -4572.7334

$ws.Range("H89").Value = 45487.39
$ws.Range("I89").Value = 126413.625
$ws.Range("J89").Value = 2326.7334
$ws.Range("K89").Value = 632068.125
$ws.Range("L89").Value = 11633.667
$ws.Range("M89").Value = -626452.125
$ws.Range("N89").Value = -22865.667

$ws.Range("H98").Value = 928.3333
$ws.Range("I98").Value = 992
$ws.Range("J98").Value = 773.7143
$ws.Range("K98").Value = 992
$ws.Range("L98").Value = 773.7143
$ws.Range("M98").Value = 506
$ws.Range("N98").Value = -3769.7143

$ws.Range("H106").Value = 1770.5
$ws.Range("I106").Value = 1745
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 1745
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -1114
$ws.Range("N106").Value = -3262

$ws.Range("H122").Value = 928.3333
$ws.Range("I122").Value = 992
$ws.Range("J122").Value = 773.7143
$ws.Range("K122").Value = 2976
$ws.Range("L122").Value = 2321.1429
$ws.Range("M122").Value = -526
$ws.Range("N122").Value = -7221.1429

$ws.Range("H123").Value = 32800
$ws.Range("J123").Value = 32800
$ws.Range("L123").Value = 32800
$ws.Range("N123").Value = -42600

$ws.Range("H135").Value = 215.46939
$ws.Range("I135").Value = 181.2889
$ws.Range("J135").Value = 600
$ws.Range("K135").Value = 1631.6001
$ws.Range("L135").Value = 5400
$ws.Range("M135").Value = 903.3998999999999
$ws.Range("N135").Value = -10470

$ws.Range("H137").Value = 2198.6042
$ws.Range("I137").Value = 1881.1666
$ws.Range("J137").Value = 2727.6667
$ws.Range("K137").Value = 5643.4998
$ws.Range("L137").Value = 8183.000100000001
$ws.Range("M137").Value = -3093.4998
$ws.Range("N137").Value = -13283.0001

$ws.Range("H141").Value = 3921.7937
$ws.Range("I141").Value = 1309.8723
$ws.Range("J141").Value = 11594.3125
$ws.Range("K141").Value = 3929.6169
$ws.Range("L141").Value = 34782.9375
$ws.Range("M141").Value = 1250.3831
$ws.Range("N141").Value = -45142.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8837.017
$ws.Range("I61").Value = 6593.064
$ws.Range("J61").Value = 16949.77
$ws.Range("K61").Value = 6593.064
$ws.Range("L61").Value = 16949.77
$ws.Range("M61").Value = -6381.064
$ws.Range("N61").Value = -17373.77

$ws.Range("H132").Value = 4314.551
$ws.Range("I132").Value = 1731.8214
$ws.Range("J132").Value = 7758.1904
$ws.Range("K132").Value = 5195.4642
$ws.Range("L132").Value = 23274.5712
$ws.Range("M132").Value = -2665.4642
$ws.Range("N132").Value = -28334.5712

$ws.Range("H136").Value = 8837.017
$ws.Range("I136").Value = 6593.064
$ws.Range("J136").Value = 16949.77
$ws.Range("K136").Value = 19779.192
$ws.Range("L136").Value = 50849.31
$ws.Range("M136").Value = -17229.192
$ws.Range("N136").Value = -55949.31

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 973
$ws.Range("I94").Value = 977.6667
$ws.Range("J94").Value = 962.5
$ws.Range("K94").Value = 977.6667
$ws.Range("L94").Value = 962.5
$ws.Range("M94").Value = -526.6667
$ws.Range("N94").Value = -1864.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1194.4117
$ws.Range("I16").Value = 588
$ws.Range("J16").Value = 1876.625
$ws.Range("K16").Value = 588
$ws.Range("L16").Value = 1876.625
$ws.Range("M16").Value = -301
$ws.Range("N16").Value = -2450.625

$ws.Range("H31").Value = 3954.1482
$ws.Range("I31").Value = 4452.2354
$ws.Range("J31").Value = 3107.4
$ws.Range("K31").Value = 4452.2354
$ws.Range("L31").Value = 3107.4
$ws.Range("M31").Value = -4157.2354
$ws.Range("N31").Value = -3697.4

$ws.Range("H34").Value = 3954.1482
$ws.Range("I34").Value = 4452.2354
$ws.Range("J34").Value = 3107.4
$ws.Range("K34").Value = 4452.2354
$ws.Range("L34").Value = 3107.4
$ws.Range("M34").Value = -4250.2354
$ws.Range("N34").Value = -3511.4

$ws.Range("H113").Value = 1194.4117
$ws.Range("I113").Value = 588
$ws.Range("J113").Value = 1876.625
$ws.Range("K113").Value = 588
$ws.Range("L113").Value = 1876.625
$ws.Range("M113").Value = 1582
$ws.Range("N113").Value = -6216.625

$ws.Range("H132").Value = 8641.143
$ws.Range("I132").Value = 13246.6
$ws.Range("J132").Value = 4454.364
$ws.Range("K132").Value = 39739.8
$ws.Range("L132").Value = 13363.092
$ws.Range("M132").Value = -37209.8
$ws.Range("N132").Value = -18423.092

$ws.Range("H134").Value = 3185.5417
$ws.Range("I134").Value = 2383.9
$ws.Range("J134").Value = 3758.1428
$ws.Range("K134").Value = 7151.700000000001
$ws.Range("L134").Value = 11274.4284
$ws.Range("M134").Value = -4616.700000000001
$ws.Range("N134").Value = -16344.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 829.8387
$ws.Range("I122").Value = 596.7692
$ws.Range("J122").Value = 998.1667
$ws.Range("K122").Value = 5370.922799999999
$ws.Range("L122").Value = 8983.5003
$ws.Range("M122").Value = -2920.922799999999
$ws.Range("N122").Value = -13883.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 42542.5
$ws.Range("J62").Value = 42542.5
$ws.Range("L62").Value = 42542.5
$ws.Range("N62").Value = -43914.5

$ws.Range("H65").Value = 42542.5
$ws.Range("J65").Value = 42542.5
$ws.Range("L65").Value = 127627.5
$ws.Range("N65").Value = -134491.5

$ws.Range("H97").Value = 1612.5
$ws.Range("I97").Value = 1520
$ws.Range("J97").Value = 1766.6666
$ws.Range("K97").Value = 1520
$ws.Range("L97").Value = 1766.6666
$ws.Range("M97").Value = -1024
$ws.Range("N97").Value = -2758.6666

$ws.Range("H102").Value = 3639.5454
$ws.Range("I102").Value = 3135.647
$ws.Range("J102").Value = 5352.8
$ws.Range("K102").Value = 3135.647
$ws.Range("L102").Value = 5352.8
$ws.Range("M102").Value = -1513.647
$ws.Range("N102").Value = -8596.799999999999

$ws.Range("H107").Value = 287.0909
$ws.Range("I107").Value = 201
$ws.Range("K107").Value = 201
$ws.Range("M107").Value = 1719

$ws.Range("H126").Value = 2337.182
$ws.Range("I126").Value = 1484.5385
$ws.Range("J126").Value = 2891.4
$ws.Range("K126").Value = 4453.6155
$ws.Range("L126").Value = 8674.200000000001
$ws.Range("M126").Value = -1983.6155
$ws.Range("N126").Value = -13614.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 25747.445
$ws.Range("J76").Value = 25747.445
$ws.Range("L76").Value = 25747.445
$ws.Range("N76").Value = -26423.445

$ws.Range("H79").Value = 25747.445
$ws.Range("J79").Value = 25747.445
$ws.Range("L79").Value = 25747.445
$ws.Range("N79").Value = -28087.445

$ws.Range("H132").Value = 5840.6772
$ws.Range("I132").Value = 5710.9585
$ws.Range("J132").Value = 6285.4287
$ws.Range("K132").Value = 17132.8755
$ws.Range("L132").Value = 18856.2861
$ws.Range("M132").Value = -14602.8755
$ws.Range("N132").Value = -23916.2861

$ws.Range("H136").Value = 2938.4146
$ws.Range("I136").Value = 1608.9365
$ws.Range("J136").Value = 7346.684
$ws.Range("K136").Value = 4826.8095
$ws.Range("L136").Value = 22040.052
$ws.Range("M136").Value = -2276.8095
$ws.Range("N136").Value = -27140.052

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2003.25
$ws.Range("I96").Value = 1456.5
$ws.Range("J96").Value = 2550
$ws.Range("K96").Value = 1456.5
$ws.Range("L96").Value = 2550
$ws.Range("M96").Value = -83.5
$ws.Range("N96").Value = -5296

$ws.Range("H132").Value = 1556.7959
$ws.Range("I132").Value = 704.09375
$ws.Range("J132").Value = 3161.8823
$ws.Range("K132").Value = 2112.28125
$ws.Range("L132").Value = 9485.6469
$ws.Range("M132").Value = 417.71875
$ws.Range("N132").Value = -14545.6469
